$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Decompose the B19 change-matrix formula into its summand parts (D19:F19)
# and running totals (D20:E20), mirroring the pattern already used on
# sheet 1 (rows 22/23).
$ws2.Range("D19").Formula = "=D3*C6*SIN(E7-G3-E6)"
$ws2.Range("E19").Formula = "=-2*E3*C7*SIN(H3)"
$ws2.Range("F19").Formula = "=-C11*SIN(E7-E11)"
$ws2.Range("D20").Formula = "=D19+E19"
$ws2.Range("E20").Formula = "=D20+F19"

# Apply a 4-decimal number format to the key result cells.
$ws2.Range("A18").NumberFormat = "0.0000"
$ws2.Range("B19").NumberFormat = "0.0000"
$ws2.Range("D20").NumberFormat = "0.0000"
$ws2.Range("E20").NumberFormat = "0.0000"

# Move the selection on sheet 2 and make sheet 1 the active tab again.
[void]$ws2.Range("E21").Select()
[void]$ws1.Activate()
